$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.659.77"
$ws.Range("E2").Value = "  +1.79%  "

$ws.Range("D3").Value = "'3.187.30"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'534.65"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").Value = "'144.27"
$ws.Range("E6").Value = "  +3.03%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  +2.21%  "

$ws.Range("D9").Value = "'7.32"
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("E10").Value = "  +1.82%  "

$ws.Range("D11").Value = "'0.427"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").Value = "'3.739.63"

$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("D14").Value = "'25.88"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").Value = "'59.711.47"
$ws.Range("E16").Value = "  +1.82%  "

$ws.Range("D17").Value = "'3.176.28"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").Value = "'13.10"
$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").Value = "'366.43"
$ws.Range("E21").Value = "  -1.44%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "'0.519"
$ws.Range("E23").Value = "  +1.34%  "

$ws.Range("D24").Value = "'69.54"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").Value = "'8.75"
$ws.Range("E25").Value = "  +9.65%  "

$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "'0.0₃0880"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").Value = "'22.22"
$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("D30").Value = "'1.89"
$ws.Range("E30").Value = "  +0.21%  "

$ws.Range("E31").Value = "  -1.15%  "

$ws.Range("D32").Value = "'5.27"
$ws.Range("E32").Value = "  +1.69%  "

$ws.Range("D33").Value = "'1.18"
$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("D34").Value = "'6.56"
$ws.Range("E34").Value = "  +4.68%  "

$ws.Range("D35").Value = "'156.82"
$ws.Range("E35").Value = "  -1.89%  "

$ws.Range("E36").Value = "  -2.32%  "

$ws.Range("D37").Value = "'2.785.48"
$ws.Range("E37").Value = "  +4.86%  "

$ws.Range("D38").Value = "'25.70"
$ws.Range("E38").Value = "  +2.00%  "

$ws.Range("D39").Value = "'0.0696"
$ws.Range("E39").Value = "  +1.53%  "

$ws.Range("E40").Value = "  -0.87%  "

$ws.Range("E41").Value = "  +0.60%  "

$ws.Range("D42").Value = "'0.0292"
$ws.Range("E42").Value = "  +2.41%  "

$ws.Range("D43").Value = "'39.24"
$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.104"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("B46").Value = "RenzoRestakedETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D46").Value = "'3.228.58"
$ws.Range("E46").Value = "  +1.17%  "

$ws.Range("D47").Value = "'0.978"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").Value = "'0.805"
$ws.Range("E48").Value = "  +5.75%  "

$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").Value = "'20.33"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("E51").Value = "  +0.01%  "

